# "automatically adding a column"
# - fixes a couple of price values
# - converts two "$x.xx" text prices back into real numbers
# - appends two new transaction rows
# - re-applies a plain 2-decimal number format across the price column
#   (replacing the old custom "$" currency format)
# - leaves a couple of pre-formatted blank cells below the table
# - tweaks the header/footer text and the current selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix / replace price values in the existing rows -----------------
$ws.Range("C3").Value = 6.65
$ws.Range("C4").Value = 7.94

# C5 / C6 used to hold text like "$9.98" - replace with real numbers
$ws.Range("C5").Value = 3.45
$ws.Range("C6").Value = 4.67

# --- append new transactions -----------------------------------------
$ws.Range("A7").Value = 1006
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 8.87

$ws.Range("A8").Value = 1007
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 6.76

# --- reformat the whole price column (header included) as 0.00 -------
$ws.Range("C1:C8").NumberFormat = "0.00"

# a couple of blank, pre-formatted cells left below the table
$ws.Range("C11:C12").NumberFormat = "0.00"
$ws.Range("C11:C12").HorizontalAlignment = -4152

# --- header / footer: add a leading space before the variable part ---
$ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12 &A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12 Page &P'

# --- selection / scroll position --------------------------------------
$ws.Range("C11").Select()
